$wb = $excel.ActiveWorkbook

# "model" sheet: add new field row for_individual (type = string)
$modelSheet = $wb.Worksheets.Item("model")
$modelSheet.Range("A15:B15").Copy() | Out-Null
$modelSheet.Range("A16:B16").PasteSpecial(-4122) | Out-Null
$modelSheet.Range("A16").Value = "string"
$modelSheet.Range("B16").Value = "for_individual"

# "table_specific_translations" sheet: add translation row for for_individual
$translationsSheet = $wb.Worksheets.Item("table_specific_translations")
$translationsSheet.Range("A16").Value = "for_individual"
$translationsSheet.Range("B16").Value = "Generate Entitlement for Individual or Benefiiciary_Entity"
